# "Einstellungen" zu "Settings" geändert
#
# Slide 23, content placeholder, the bullet "Projekt <arrow> Einstellungen
# <arrow> CI/CD <arrow> Runners" gets its "Einstellungen" run split so a new
# "<arrow> Settings " run is inserted before the remaining "<arrow> CI/CD
# <arrow> Runners" run.

$p = $ppt.ActivePresentation
$arrow = [char]0xF0E0
$needle = $arrow + " Einstellungen "
$replacement = $arrow + " Settings "

$found = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }

        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }

        $tr = $tf.TextRange
        $paraCount = $tr.Paragraphs().Count

        for ($pi = 1; $pi -le $paraCount; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            $text = $para.Text

            # 0-based index of the "<arrow> Einstellungen " run inside the
            # paragraph text.
            $idx = $text.IndexOf($needle)
            if ($idx -lt 0) { continue }

            # TextRange.Characters is 1-based.
            $range = $para.Characters($idx + 1, $needle.Length)
            $range.Text = $replacement

            $found = $true
        }
    }
}

if (-not $found) {
    throw "Could not find the 'Einstellungen' run to update"
}
